# Generate Report for Handoff
#
# The localization-status report is regenerated: the 9f969ab8 markdown
# file has a fresh handoff xliff, so its "Latest HO Xliff Generate Date"
# (Overview sheet) and "Latest Handoff Datetime" (per-locale sheets) move
# forward to the new timestamps produced by this handoff run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G6").Value = "2016-08-23 22:39:33"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H6").Value = "2016-08-23 22:39:28"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H6").Value = "2016-08-23 22:39:33"
